{"js": "// Replace each old arithmetic expression with its updated counterpart.\n// The mapping below is derived from the target diff; every \"old\" value is\n// unique in the document, so sequential body.search() + insertText(Replace)\n// calls (no wildcards) are safe and order-independent, and they preserve the\n// run/paragraph formatting of the matched text.\nconst pairs = [\n  ['76-26=50', '54+26=80'],\n  ['83-14=69', '14+61=75'],\n  ['95-63=32', '21+52=73'],\n  ['62-58=4', '58-47=11'],\n  ['56-45=11', '80-66=14'],\n  ['51+29=80', '6+69=75'],\n  ['47-44=3', '43+39=82'],\n  ['66-10=56', '68+7=75'],\n  ['95-52=43', '95-78=17'],\n  ['82-20=62', '1+16=17'],\n  ['46-33=13', '88-79=9'],\n  ['9+31=40', '46+40=86'],\n  ['76-71=5', '80+16=96'],\n  ['32+6=38', '53+46=99'],\n  ['94-42=52', '38+23=61'],\n  ['64-29=35', '42-17=25'],\n  ['18-5=13', '37+21=58'],\n  ['79-60=19', '78-57=21'],\n  ['7+77=84', '89-0=89'],\n  ['54+19=73', '98+1=99'],\n  ['48-18=30', '9-1=8'],\n  ['70-9=61', '74-62=12'],\n  ['39+18=57', '8+53=61'],\n  ['83-56=27', '87-34=53'],\n  ['38-27=11', '18+32=50'],\n  ['74-51=23', '58-37=21'],\n  ['28-7=21', '83-16=67'],\n  ['69+22=91', '58-53=5'],\n  ['71-30=41', '99-48=51'],\n  ['39+33=72', '49+34=83'],\n  ['13+26=39', '98-14=84'],\n  ['4+62=66', '25+24=49'],\n  ['20-0=20', '33+13=46'],\n  ['59+2=61', '57+17=74'],\n  ['92-85=7', '66-26=40'],\n  ['34-12=22', '22-16=6'],\n  ['56+0=56', '66+24=90'],\n  ['93-62=31', '7+66=73'],\n  ['27+7=34', '28+17=45'],\n  ['65+27=92', '83-70=13'],\n  ['81-3=78', '18+65=83'],\n  ['16+58=74', '89-54=35'],\n  ['23-5=18', '3+88=91'],\n  ['39+55=94', '11+72=83'],\n  ['85-45=40', '45-25=20'],\n  ['9+87=96', '41+12=53'],\n  ['61+0=61', '24+27=51'],\n  ['38-15=23', '95-74=21'],\n  ['52+20=72', '6+47=53'],\n  ['32-7=25', '91-62=29'],\n  ['91-87=4', '39+43=82'],\n  ['22+40=62', '60+33=93'],\n  ['91-0=91', '0+83=83'],\n  ['73-18=55', '71+13=84'],\n  ['75+3=78', '50+15=65'],\n  ['81-56=25', '71-26=45'],\n  ['51-46=5', '37-25=12'],\n  ['8+58=66', '96-8=88'],\n  ['62-54=8', '27-14=13'],\n  ['94-25=69', '40+23=63'],\n  ['42+54=96', '64+1=65'],\n  ['47+45=92', '8+68=76'],\n  ['85-22=63', '29+59=88'],\n  ['83-2=81', '53-42=11'],\n  ['9+81=90', '59+3=62'],\n  ['68-52=16', '32+37=69'],\n  ['48-9=39', '66-21=45'],\n  ['76-25=51', '71-46=25'],\n  ['4+54=58', '13-9=4'],\n  ['17-8=9', '17+72=89'],\n  ['58+17=75', '43-27=16'],\n  ['53+11=64', '57+26=83'],\n  ['17+26=43', '17+33=50'],\n  ['98-68=30', '9-7=2'],\n  ['76-35=41', '79+9=88'],\n  ['24+29=53', '23+21=44'],\n  ['78-9=69', '66-64=2'],\n  ['76+22=98', '88-86=2'],\n  ['27+49=76', '93-10=83'],\n  ['15+23=38', '58+6=64'],\n  ['64-15=49', '25+63=88'],\n  ['84-10=74', '36-10=26'],\n  ['16+50=66', '5+26=31'],\n  ['32+38=70', '62+24=86'],\n  ['60-35=25', '44+11=55'],\n  ['7+1=8', '45+3=48'],\n  ['66+30=96', '26+66=92'],\n  ['56-17=39', '29-26=3'],\n  ['71-34=37', '15+68=83'],\n  ['21+32=53', '78+20=98'],\n  ['46+15=61', '5+61=66'],\n  ['64+0=64', '54-5=49'],\n  ['12+45=57', '64+11=75'],\n  ['99-67=32', '83+0=83'],\n  ['35+8=43', '87+6=93'],\n  ['79+11=90', '99-28=71'],\n  ['54+10=64', '34-0=34'],\n  ['30+7=37', '70-22=48'],\n  ['49+7=56', '47-4=43'],\n  ['10+82=92', '19+62=81'],\n];\n\nconst body = context.document.body;\nfor (const [oldText, newText] of pairs) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (const r of results.items) {\n    r.insertText(newText, Word.InsertLocation.replace);\n  }\n}\nawait context.sync();\n", "ps1": "# Replace each old arithmetic expression with its updated counterpart.\n# The mapping below is derived from the target diff; every \"old\" value is\n# unique in the document, so a plain sequential Find/Replace (no wildcards)\n# against the whole document body is safe and order-independent.\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @('76-26=50', '54+26=80'),\n    @('83-14=69', '14+61=75'),\n    @('95-63=32', '21+52=73'),\n    @('62-58=4', '58-47=11'),\n    @('56-45=11', '80-66=14'),\n    @('51+29=80', '6+69=75'),\n    @('47-44=3', '43+39=82'),\n    @('66-10=56', '68+7=75'),\n    @('95-52=43', '95-78=17'),\n    @('82-20=62', '1+16=17'),\n    @('46-33=13', '88-79=9'),\n    @('9+31=40', '46+40=86'),\n    @('76-71=5', '80+16=96'),\n    @('32+6=38', '53+46=99'),\n    @('94-42=52', '38+23=61'),\n    @('64-29=35', '42-17=25'),\n    @('18-5=13', '37+21=58'),\n    @('79-60=19', '78-57=21'),\n    @('7+77=84', '89-0=89'),\n    @('54+19=73', '98+1=99'),\n    @('48-18=30', '9-1=8'),\n    @('70-9=61', '74-62=12'),\n    @('39+18=57', '8+53=61'),\n    @('83-56=27', '87-34=53'),\n    @('38-27=11', '18+32=50'),\n    @('74-51=23', '58-37=21'),\n    @('28-7=21', '83-16=67'),\n    @('69+22=91', '58-53=5'),\n    @('71-30=41', '99-48=51'),\n    @('39+33=72', '49+34=83'),\n    @('13+26=39', '98-14=84'),\n    @('4+62=66', '25+24=49'),\n    @('20-0=20', '33+13=46'),\n    @('59+2=61', '57+17=74'),\n    @('92-85=7', '66-26=40'),\n    @('34-12=22', '22-16=6'),\n    @('56+0=56', '66+24=90'),\n    @('93-62=31', '7+66=73'),\n    @('27+7=34', '28+17=45'),\n    @('65+27=92', '83-70=13'),\n    @('81-3=78', '18+65=83'),\n    @('16+58=74', '89-54=35'),\n    @('23-5=18', '3+88=91'),\n    @('39+55=94', '11+72=83'),\n    @('85-45=40', '45-25=20'),\n    @('9+87=96', '41+12=53'),\n    @('61+0=61', '24+27=51'),\n    @('38-15=23', '95-74=21'),\n    @('52+20=72', '6+47=53'),\n    @('32-7=25', '91-62=29'),\n    @('91-87=4', '39+43=82'),\n    @('22+40=62', '60+33=93'),\n    @('91-0=91', '0+83=83'),\n    @('73-18=55', '71+13=84'),\n    @('75+3=78', '50+15=65'),\n    @('81-56=25', '71-26=45'),\n    @('51-46=5', '37-25=12'),\n    @('8+58=66', '96-8=88'),\n    @('62-54=8', '27-14=13'),\n    @('94-25=69', '40+23=63'),\n    @('42+54=96', '64+1=65'),\n    @('47+45=92', '8+68=76'),\n    @('85-22=63', '29+59=88'),\n    @('83-2=81', '53-42=11'),\n    @('9+81=90', '59+3=62'),\n    @('68-52=16', '32+37=69'),\n    @('48-9=39', '66-21=45'),\n    @('76-25=51', '71-46=25'),\n    @('4+54=58', '13-9=4'),\n    @('17-8=9', '17+72=89'),\n    @('58+17=75', '43-27=16'),\n    @('53+11=64', '57+26=83'),\n    @('17+26=43', '17+33=50'),\n    @('98-68=30', '9-7=2'),\n    @('76-35=41', '79+9=88'),\n    @('24+29=53', '23+21=44'),\n    @('78-9=69', '66-64=2'),\n    @('76+22=98', '88-86=2'),\n    @('27+49=76', '93-10=83'),\n    @('15+23=38', '58+6=64'),\n    @('64-15=49', '25+63=88'),\n    @('84-10=74', '36-10=26'),\n    @('16+50=66', '5+26=31'),\n    @('32+38=70', '62+24=86'),\n    @('60-35=25', '44+11=55'),\n    @('7+1=8', '45+3=48'),\n    @('66+30=96', '26+66=92'),\n    @('56-17=39', '29-26=3'),\n    @('71-34=37', '15+68=83'),\n    @('21+32=53', '78+20=98'),\n    @('46+15=61', '5+61=66'),\n    @('64+0=64', '54-5=49'),\n    @('12+45=57', '64+11=75'),\n    @('99-67=32', '83+0=83'),\n    @('35+8=43', '87+6=93'),\n    @('79+11=90', '99-28=71'),\n    @('54+10=64', '34-0=34'),\n    @('30+7=37', '70-22=48'),\n    @('49+7=56', '47-4=43'),\n    @('10+82=92', '19+62=81'),\n)\n\nforeach ($pair in $pairs) {\n    $old = $pair[0]\n    $new = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $old\n    $find.Replacement.Text = $new\n    $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n}\n"}
